$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying NATMI computation was re-run against the updated TPM matrix.
# Only the "Target cluster = FAPs" pairs survive the re-run (the old "Target
# cluster = ECs" rows are gone), so the 8 original data rows collapse down to
# 4 - one per Sending cluster - with freshly recomputed specificity scores.
# Rows 2-5 below hold the new values (old rows 3/5/7/9 are the closest
# analogues, reused as the "Sending cluster" order), and the former rows 6-9
# are cleared since the table now only spans A1:T5.

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ntn1"
$ws.Range("C2").Value = "Dcc"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.629231666666667
$ws.Range("H2").Value = 7.887695
$ws.Range("I2").Value = 0.1414315557047068
$ws.Range("J2").Value = 0.1414315557047067
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.05453333333333333
$ws.Range("N2").Value = 0.1636
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.1433807668888889
$ws.Range("R2").Value = 1.290426902
$ws.Range("S2").Value = 0.1414315557047068
$ws.Range("T2").Value = 0.1414315557047067

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Ntn1"
$ws.Range("C3").Value = "Dcc"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 11.42765333333333
$ws.Range("H3").Value = 34.28296
$ws.Range("I3").Value = 0.6147160060020365
$ws.Range("J3").Value = 0.6147160060020365
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.05453333333333333
$ws.Range("N3").Value = 0.1636
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.6231880284444444
$ws.Range("R3").Value = 5.608692256
$ws.Range("S3").Value = 0.6147160060020365
$ws.Range("T3").Value = 0.6147160060020365

# Row 4
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Ntn1"
$ws.Range("C4").Value = "Dcc"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4.24731
$ws.Range("H4").Value = 12.74193
$ws.Range("I4").Value = 0.2284711798035388
$ws.Range("J4").Value = 0.2284711798035388
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.05453333333333333
$ws.Range("N4").Value = 0.1636
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.231619972
$ws.Range("R4").Value = 2.084579748
$ws.Range("S4").Value = 0.2284711798035388
$ws.Range("T4").Value = 0.2284711798035388

# Row 5
$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("B5").Value = "Ntn1"
$ws.Range("C5").Value = "Dcc"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.2859396666666667
$ws.Range("H5").Value = 0.8578190000000001
$ws.Range("I5").Value = 0.01538125848971795
$ws.Range("J5").Value = 0.01538125848971795
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.05453333333333333
$ws.Range("N5").Value = 0.1636
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 0.01559324315555556
$ws.Range("R5").Value = 0.1403391884
$ws.Range("S5").Value = 0.01538125848971795
$ws.Range("T5").Value = 0.01538125848971795

# Remove the now-obsolete rows (old MuSCs/Resolving-Mac "ECs-target" & duplicate rows
# have been folded into rows 2-5 above); clear what used to be rows 6-9 so the sheet
# shrinks back down to a 4-data-row table (dimension A1:T5).
$ws.Range("A6:T9").ClearContents()

